$wb = $excel.ActiveWorkbook
$log = $wb.Worksheets.Item("Log")

# Fill in the missing data for the "15.Mix" test rows (26-28) on the Log sheet.
$log.Range("C26").Value = 0.8125
$log.Range("G26").Value = 42538
$log.Range("H26").Value = "12.57.txt"

$log.Range("C27").Value = 0.70703125
$log.Range("G27").Value = 42538
$log.Range("H27").Value = "13.12.txt"

$log.Range("C28").Value = 0.8125
$log.Range("G28").Value = 42538
$log.Range("H28").Value = "14.20.txt"

$new = $wb.Worksheets.Add($null, $log)
$new.Name = "15-Series"

# copy header row
$log.Range("A1:I1").Copy()
$new.Range("A1").PasteSpecial()

# copy the 15.* rows (20-28) into the new sheet starting row 2
$log.Range("A20:H28").Copy()
$new.Range("A2").PasteSpecial()

Write-Host "done"
